$d = $word.ActiveDocument

# The last paragraph currently reads:
#   "NEXT: Complete Notifications class and the tournament class."
# and carries the _GoBack bookmark at its end. We append five new
# paragraphs after it (mirroring the existing "date header" + bullet
# pattern used throughout the log), and the _GoBack bookmark will then
# sit at the very end of the new content, matching Word's normal
# behavior of the bookmark tracking the last edited location.

$last = $d.Paragraphs.Last
$r = $last.Range
$r.InsertParagraphAfter()

# 1) Blank spacer paragraph (no bullet numbering), like the blank line
#    separating each day's entries elsewhere in the log.
$p1 = $d.Paragraphs.Last
$p1.Range.ListFormat.RemoveNumbers()
$p1.LeftIndent = 18
$p1.FirstLineIndent = -18

$p1.Range.InsertParagraphAfter()

# 2) Bold date header: "09/20/2016 (1hr)"
$p2 = $d.Paragraphs.Last
$p2.Range.ListFormat.RemoveNumbers()
$p2.LeftIndent = 18
$p2.FirstLineIndent = -18
$r2 = $p2.Range
$r2.InsertAfter("09/20/2016 (1hr)")
$r2.Bold = 1

$d.Paragraphs.Last.Range.InsertParagraphAfter()

# 3) Bullet: Notifications class note
$p3 = $d.Paragraphs.Last
$r3 = $p3.Range
$r3.Collapse(0)
$r3.InsertAfter("Almost fully implemented the Notifications class for the classes written so far.")

$d.Paragraphs.Last.Range.InsertParagraphAfter()

# 4) Bullet: validation ordering note
$p4 = $d.Paragraphs.Last
$r4 = $p4.Range
$r4.Collapse(0)
$r4.InsertAfter("Potential issue is that the input validation waits until all the coordinates and path choices have been entered. Nothing fails, it" + [char]0x2019 + "s just weird to let user select path when the destination is out of bounds already.")

$d.Paragraphs.Last.Range.InsertParagraphAfter()

# 5) Bullet: board-printing note (keeps the _GoBack bookmark at the end)
$p5 = $d.Paragraphs.Last
$r5 = $p5.Range
$r5.Collapse(0)
$r5.InsertAfter("Modified code to printing board only after each user makes a valid move")

Write-Output "done"
